# Scheduled market-data refresh: update currentAveragePrice/LevePrice/LeveProfit
# columns (H:N) on the per-job sheets. Values come from the latest market
# snapshot; a handful of rows also gain/lose an NQ or HQ profit cell because
# the refreshed recipe no longer has (or now has) an HQ variant.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3943.639
$ws.Range("I15").Value = 3943.639
$ws.Range("K15").Value = 11830.917
$ws.Range("M15").Value = -11661.917

$ws.Range("H18").Value = 1107.6666
$ws.Range("J18").Value = 1749.75
$ws.Range("L18").Value = 1749.75
$ws.Range("N18").Value = -2317.75

$ws.Range("H86").Value = 1003225.94
$ws.Range("I86").Value = 1335767.9
$ws.Range("K86").Value = 1335767.9
$ws.Range("M86").Value = -1334644.9

$ws.Range("H89").Value = 1003225.94
$ws.Range("I89").Value = 1335767.9
$ws.Range("K89").Value = 6678839.5
$ws.Range("M89").Value = -6673223.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 164.05882
$ws.Range("J5").Value = 257.85715
$ws.Range("L5").Value = 257.85715
$ws.Range("N5").Value = -481.85715

$ws.Range("H45").Value = 22729628
$ws.Range("J45").Value = 4040
$ws.Range("L45").Value = 4040
$ws.Range("N45").Value = -4794

$ws.Range("H74").Value = 10405140
$ws.Range("I74").Value = 20837062
$ws.Range("K74").Value = 20837062
$ws.Range("M74").Value = -20836188

$ws.Range("H77").Value = 10405140
$ws.Range("I77").Value = 20837062
$ws.Range("K77").Value = 104185310
$ws.Range("M77").Value = -104180942

$ws.Range("H88").Value = 2900.2593
$ws.Range("I88").Value = 2247
$ws.Range("J88").Value = 3013.8696
$ws.Range("K88").Value = 2247
$ws.Range("L88").Value = 3013.8696
$ws.Range("M88").Value = -1841
$ws.Range("N88").Value = -3825.8696

$ws.Range("H91").Value = 2900.2593
$ws.Range("I91").Value = 2247
$ws.Range("J91").Value = 3013.8696
$ws.Range("K91").Value = 2247
$ws.Range("L91").Value = 3013.8696
$ws.Range("M91").Value = -843
$ws.Range("N91").Value = -5821.8696

$ws.Range("H97").Value = 2014.95
$ws.Range("I97").Value = 1786.3334
$ws.Range("K97").Value = 1786.3334
$ws.Range("M97").Value = -1290.3334

$ws.Range("H122").Value = 2902
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 2202.6667
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 6608.000100000001
$ws.Range("M122").Value = -12550
$ws.Range("N122").Value = -11508.0001

$ws.Range("H132").Value = 2348.9658
$ws.Range("I132").Value = 1441.7715
$ws.Range("K132").Value = 4325.3145
$ws.Range("M132").Value = -1795.3145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 164.05882
$ws.Range("J4").Value = 257.85715
$ws.Range("L4").Value = 257.85715
$ws.Range("N4").Value = -487.85715

$ws.Range("H7").Value = 2438.6667
$ws.Range("I7").Value = 926.2
$ws.Range("K7").Value = 926.2
$ws.Range("M7").Value = -813.2

$ws.Range("H20").Value = 9003.875
$ws.Range("I20").Value = 9072.666999999999
$ws.Range("J20").Value = 8797.5
$ws.Range("K20").Value = 9072.666999999999
$ws.Range("L20").Value = 8797.5
$ws.Range("M20").Value = -8825.666999999999
$ws.Range("N20").Value = -9291.5

$ws.Range("H26").Value = 35788
$ws.Range("I26").Value = 17980
$ws.Range("J26").Value = 62500
$ws.Range("K26").Value = 17980
$ws.Range("L26").Value = 62500
$ws.Range("M26").Value = -17688
$ws.Range("N26").Value = -63084

$ws.Range("H40").Value = 234949.5
$ws.Range("J40").Value = 234949.5
$ws.Range("L40").Value = 234949.5
$ws.Range("N40").Value = -235479.5

$ws.Range("H94").Value = 1451.7391
$ws.Range("I94").Value = 1438.8667
$ws.Range("J94").Value = 1475.875
$ws.Range("K94").Value = 1438.8667
$ws.Range("L94").Value = 1475.875
$ws.Range("M94").Value = -987.8667
$ws.Range("N94").Value = -2377.875

$ws.Range("H95").Value = 34999.5
$ws.Range("J95").Value = 34999.5
$ws.Range("L95").Value = 34999.5
$ws.Range("N95").Value = -40491.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 107.5
$ws.Range("I7").Value = 147.14285
$ws.Range("K7").Value = 147.14285
$ws.Range("M7").Value = -34.14285000000001

$ws.Range("H132").Value = 1506.5625
$ws.Range("I132").Value = 1578.9286
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 4736.7858
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -2206.7858
$ws.Range("N132").Value = -8060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4986.448
$ws.Range("I131").Value = 4200.357
$ws.Range("K131").Value = 12601.071
$ws.Range("M131").Value = -7561.071

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 51673.4
$ws.Range("J24").Value = 53364.355
$ws.Range("L24").Value = 53364.355
$ws.Range("N24").Value = -53710.355

$ws.Range("H93").Value = 59990
$ws.Range("J93").Value = 59990
$ws.Range("L93").Value = 59990
$ws.Range("N93").Value = -63734

$ws.Range("H122").Value = 2260.3
$ws.Range("I122").Value = 1934.3334
$ws.Range("J122").Value = 2749.25
$ws.Range("K122").Value = 5803.0002
$ws.Range("L122").Value = 8247.75
$ws.Range("M122").Value = -3353.0002
$ws.Range("N122").Value = -13147.75

$ws.Range("H132").Value = 25005990
$ws.Range("I132").Value = 30306488
$ws.Range("J132").Value = 17938.285
$ws.Range("K132").Value = 90919464
$ws.Range("L132").Value = 53814.855
$ws.Range("M132").Value = -90916934
$ws.Range("N132").Value = -58874.855

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2186.5715
$ws.Range("I82").Value = 500.5
$ws.Range("J82").Value = 2861
$ws.Range("K82").Value = 500.5
$ws.Range("L82").Value = 2861
$ws.Range("M82").Value = -139.5
$ws.Range("N82").Value = -3583

$ws.Range("H85").Value = 2186.5715
$ws.Range("I85").Value = 500.5
$ws.Range("J85").Value = 2861
$ws.Range("K85").Value = 500.5
$ws.Range("L85").Value = 2861
$ws.Range("M85").Value = 747.5
$ws.Range("N85").Value = -5357

$ws.Range("H122").Value = 5494.684
$ws.Range("I122").Value = 4826.6
$ws.Range("K122").Value = 14479.8
$ws.Range("M122").Value = -12029.8

$ws.Range("H136").Value = 167428.62
$ws.Range("I136").Value = 36741.332
$ws.Range("K136").Value = 110223.996
$ws.Range("M136").Value = -107673.996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 18825082
$ws.Range("J75").Value = 18825082
$ws.Range("L75").Value = 18825082
$ws.Range("N75").Value = -18826954

$ws.Range("H78").Value = 18825082
$ws.Range("J78").Value = 18825082
$ws.Range("L78").Value = 56475246
$ws.Range("N78").Value = -56484606

$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()

$ws.Range("H126").Value = 3164.8
$ws.Range("I126").Value = 3475
$ws.Range("J126").Value = 2699.5
$ws.Range("K126").Value = 10425
$ws.Range("L126").Value = 8098.5
$ws.Range("M126").Value = -7955
$ws.Range("N126").Value = -13038.5

$ws.Range("H132").Value = 2142.875
$ws.Range("I132").Value = 1619.0667
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 4857.2001
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -2327.2001
$ws.Range("N132").Value = -35060

$ws.Range("H136").Value = 4720.684
$ws.Range("J136").Value = 6000
$ws.Range("L136").Value = 18000
$ws.Range("N136").Value = -23100
